$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("network")

# New header values for B1:BA1 -- these mirror the row labels in column A
# (genes controlling), replacing the old all-caps duplicate shared strings
# that used to live only in this header row.
$values = @(
    "ADH5","AGP2","AIM3","AIM4","ALG1","AMN1","APD1","ARA1","ARL1","ATG14",
    "BEM1","BMT2","CBP6","CCZ1","CDC28","CKS1","CMD1","CNS1","COS111","CYC8",
    "DER1","DTR1","DUR1","ECM31","EHT1","EXO5","EXO84","FES1","FLO9","FTH1",
    "FZO1","GDT1","GRS1","HSL7","ICS2","IFA38","IML3","IRA1","KTR3","KTR4",
    "LDH1","LYS2","MAK5","MBA1","MCM7","MEC1","MED8","MMS4","MRPL36","MRPS9",
    "MSI1","MUD1"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 2
    $ws.Cells.Item(1, $col).Value = $values[$i]
}

# Update the sheet's selection to match the new state (B1:BA1 selected, active cell B1)
$ws.Range("B1:BA1").Select()
